$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values keep their original text representation
# (values like "309.52" would otherwise be auto-converted to numbers by Excel)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.858.81'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '1.809.98'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '309.52'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D7').Value = '0.4643'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('D8').Value = '0.3700'
$ws.Range('E8').Value = '  -2.45%  '
$ws.Range('D9').Value = '0.07368'
$ws.Range('E9').Value = '  -0.69%  '
$ws.Range('D10').Value = '0.8754'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = '20.44'
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').Value = '1.802.55'
$ws.Range('E12').Value = '  -1.34%  '
$ws.Range('D13').Value = '5.362'
$ws.Range('E13').Value = '  -1.44%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '0.07074'
$ws.Range('E14').Value = '  -0.54%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '6.508'
$ws.Range('E15').Value = '  -3.16%  '
$ws.Range('D16').Value = '91.55'
$ws.Range('E16').Value = '  -1.83%  '
$ws.Range('E17').Value = '  +0.17%  '
$ws.Range('D18').Value = '0.000008708'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('E20').Value = '  -1.98%  '
$ws.Range('D21').Value = '26.860.29'
$ws.Range('D22').Value = '5.317'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D24').Value = '1.983.53'
$ws.Range('E24').Value = '  -3.24%  '
$ws.Range('D25').Value = '1.899'
$ws.Range('D26').Value = '151.56'
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('E27').Value = '  -1.37%  '
$ws.Range('D28').Value = '2.152'
$ws.Range('E28').Value = '  -5.18%  '
$ws.Range('D29').Value = '5.312'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('D30').Value = '115.88'
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('E31').Value = '  -0.49%  '
$ws.Range('D32').Value = '0.7549'
$ws.Range('E32').Value = '  -5.50%  '
$ws.Range('D33').Value = '1.155'
$ws.Range('E33').Value = '  -3.62%  '
$ws.Range('D34').Value = '4.465'
$ws.Range('E34').Value = '  -1.87%  '
$ws.Range('D35').Value = '2.915'
$ws.Range('E35').Value = '  -0.61%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('E37').Value = '  -0.13%  '
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('D39').Value = '2.445'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').Value = '0.05258'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('D41').Value = '2.913'
$ws.Range('E41').Value = '  +0.70%  '
$ws.Range('D42').Value = '0.5311'
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('D43').Value = '7.176'
$ws.Range('E43').Value = '  -2.14%  '
$ws.Range('E44').Value = '  -2.75%  '
$ws.Range('D45').Value = '8.465'
$ws.Range('E45').Value = '  -2.60%  '
$ws.Range('D46').Value = '0.4946'
$ws.Range('D47').Value = '10.39'
$ws.Range('E47').Value = '  -1.73%  '
$ws.Range('E48').Value = '  -0.78%  '
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('E50').Value = '  -2.54%  '
$ws.Range('D51').Value = '0.06291'
$ws.Range('E51').Value = '  -1.56%  '
